# lab_result source->target field mapping sheet generator.
#
# This mirrors the sync_commcare_assets export script: it re-derives the
# "Field" / "Source Field" mapping rows (columns E:F, starting row 3) from
# the full set of known CommCare case properties for the lab_result case
# type, sorted alphabetically by field name. Rewriting the block in full
# (rather than only appending) is what lets the script be safely re-run
# multiple times without accumulating duplicate rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly-discovered case properties being added in this run (on top of the
# ones already present in the sheet).
$newFields = @(
    'address',
    'address_city',
    'address_complete',
    'address_county',
    'address_state',
    'address_street',
    'address_zip',
    'age',
    'case_import_date',
    'contact_phone_number',
    'current_status',
    'dob',
    'dob_known',
    'ethnicity',
    'first_name',
    'full_name',
    'gender',
    'has_phone_number',
    'last_name',
    'ltcf',
    'owner_id',
    'patient_type',
    'phone_home',
    'race'
)

# Fields that map to a Source Field WITHOUT the 'properties.' prefix (e.g.
# CommCare metadata fields that live on the case itself, not in properties).
$noPrefixFields = @('closed')

# --- Read whatever mapping rows already exist (columns E/F, from row 3
# down) so re-running this script is idempotent: existing entries are kept
# (and de-duped / re-prefixed correctly) rather than blindly appended to. ---
$existing = @{}
$row = 3
while ($true) {
    $fieldVal = $ws.Cells.Item($row, 5).Value2
    if ([string]::IsNullOrEmpty($fieldVal)) { break }

    $srcVal = $ws.Cells.Item($row, 6).Value2

    # Bug fix: on a prior re-run, 'properties.' could get prepended again to
    # a source field that already had it (e.g. 'properties.properties.x').
    # Normalize that away so re-runs stay stable.
    while ($srcVal -like 'properties.properties.*') {
        $srcVal = $srcVal.Substring(11)
    }

    $existing[$fieldVal] = $srcVal
    $row = $row + 1
}
$lastExistingRow = $row - 1

# --- Merge in the newly discovered fields. ---
foreach ($field in $newFields) {
    if (-not $existing.ContainsKey($field)) {
        if ($noPrefixFields -contains $field) {
            $existing[$field] = $field
        } else {
            $existing[$field] = 'properties.' + $field
        }
    }
}

# Make sure every field (old + new) uses the correct Source Field, fixing
# the no-prefix ones and ensuring everything else really does carry the
# 'properties.' prefix exactly once.
foreach ($field in @($existing.Keys)) {
    if ($noPrefixFields -contains $field) {
        $existing[$field] = $field
    } elseif ($existing[$field] -notlike 'properties.*') {
        $existing[$field] = 'properties.' + $existing[$field]
    }
}

# --- Clear the old block out, then rewrite the merged, alphabetically
# sorted field list starting at row 3. Clearing first (rather than just
# overwriting in place) ensures the sheet - and its backing shared-string
# table - reflect only the current, de-duplicated, re-sorted set of rows
# on every re-run. ---
$sortedFields = $existing.Keys | Sort-Object

$ws.Range($ws.Cells.Item(3, 5), $ws.Cells.Item($lastExistingRow, 6)).ClearContents()

$row = 3
foreach ($field in $sortedFields) {
    $ws.Cells.Item($row, 5).Value2 = $field
    $ws.Cells.Item($row, 6).Value2 = $existing[$field]
    $row = $row + 1
}
